$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.726.13"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "2.648.77"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.43%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("E10").Value = "  +1.73%  "

$ws.Range("E11").Value = "  +2.90%  "

$ws.Range("E12").Value = "  +3.21%  "

$ws.Range("D13").Value = "3.120.13"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.39%  "

$ws.Range("D15").Value = "60.726.07"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("D17").Value = "2.661.54"
$ws.Range("E17").Value = "  +1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.531"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.94%  "

$ws.Range("E28").Value = "  +10.67%  "

$ws.Range("D29").Value = "0.0₃0811"
$ws.Range("E29").Value = "  +3.49%  "

$ws.Range("E30").Value = "  +4.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.51%  "

$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.75%  "

$ws.Range("E36").Value = "  +8.02%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "333.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.10%  "

$ws.Range("E39").Value = "  +4.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("E41").Value = "  +3.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.616"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.59%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0560"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  +2.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "
